$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for years 2005-2009 (rows 2 through 6),
# shifting the remaining data (2010-2015) up.
$ws.Range("A2:G6").Delete(-4162)
